$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three countries leapfrogged their neighbour in the ranking once the
# "Casos totales" figures were refreshed, so the two rows swap names
# (their historic stats ride along to the row they moved into).
$ws.Range("A54").Value  = "Bielorrusia"   # was Honduras
$ws.Range("A55").Value  = "Honduras"      # was Bielorrusia

$ws.Range("A95").Value  = "Albania"       # was Noruega
$ws.Range("A96").Value  = "Noruega"       # was Albania

$ws.Range("A134").Value = "Sri Lanka"                      # was Republica de Africa Central
$ws.Range("A135").Value = "Republica de Africa Central"    # was Sri Lanka

# Refreshed Casos totales / Nuevos casos / Casos activos / Recuperados /
# Muertes hoy / Muertes counts (columns B,C,D,E,G,H) for each affected row.
# Row 4: Estados Unidos
$ws.Range("B4").Value = 8044067
$ws.Range("C4").Value = 6278
$ws.Range("D4").Value = 5196547
$ws.Range("E4").Value = 2627346
$ws.Range("G4").Value = 163
$ws.Range("H4").Value = 220174

# Row 5: India
$ws.Range("B5").Value = 7205923
$ws.Range("C5").Value = 32358
$ws.Range("D5").Value = 6255622
$ws.Range("E5").Value = 840183
$ws.Range("G5").Value = 224
$ws.Range("H5").Value = 110118

# Row 17: Chile
$ws.Range("B17").Value = 484280
$ws.Range("C17").Value = 1392
$ws.Range("D17").Value = 456499
$ws.Range("E17").Value = 14385
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 13396

# Row 25: Alemania
$ws.Range("B25").Value = 333311
$ws.Range("C25").Value = 2217
$ws.Range("E25").Value = 44482
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 9729

# Row 30: Canada
$ws.Range("B30").Value = 184392
$ws.Range("C30").Value = 1553
$ws.Range("D30").Value = 154882
$ws.Range("E30").Value = 19874
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = 9636

# Row 40: Republica Dominicana
$ws.Range("B40").Value = 119008
$ws.Range("C40").Value = 165
$ws.Range("D40").Value = 94941
$ws.Range("E40").Value = 21884
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 2183

# Row 48: Guatemala
$ws.Range("B48").Value = 98380
$ws.Range("C48").Value = 554
$ws.Range("D48").Value = 87744
$ws.Range("E48").Value = 7226
$ws.Range("G48").Value = 23
$ws.Range("H48").Value = 3410

# Row 49: Japon
$ws.Range("B49").Value = 89673
$ws.Range("C49").Value = 326
$ws.Range("D49").Value = 82621
$ws.Range("E49").Value = 5418
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 1634

# Row 54: Bielorrusia
$ws.Range("B54").Value = 84524
$ws.Range("C54").Value = 526
$ws.Range("D54").Value = 77797
$ws.Range("E54").Value = 5821
$ws.Range("G54").Value = 5
$ws.Range("H54").Value = 906

# Row 55: Honduras
$ws.Range("B55").Value = 84413
$ws.Range("C55").Value = 332
$ws.Range("D55").Value = 32214
$ws.Range("E55").Value = 49678
$ws.Range("G55").Value = 9
$ws.Range("H55").Value = 2521

# Row 59: Moldavia
$ws.Range("B59").Value = 63275
$ws.Range("C59").Value = 657
$ws.Range("D59").Value = 45102
$ws.Range("E59").Value = 16678
$ws.Range("G59").Value = 17
$ws.Range("H59").Value = 1495

# Row 87: Grecia
$ws.Range("B87").Value = 23060
$ws.Range("C87").Value = 408
$ws.Range("E87").Value = 12609
$ws.Range("G87").Value = 6
$ws.Range("H87").Value = 462

# Row 95: Albania
$ws.Range("B95").Value = 15752
$ws.Range("C95").Value = 182
$ws.Range("D95").Value = 9675
$ws.Range("E95").Value = 5648
$ws.Range("G95").Value = 5
$ws.Range("H95").Value = 429

# Row 96: Noruega
$ws.Range("B96").Value = 15730
$ws.Range("C96").Value = 91
$ws.Range("D96").Value = 11863
$ws.Range("E96").Value = 3590
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 277

# Row 99: Montenegro
$ws.Range("B99").Value = 14268
$ws.Range("C99").Value = 218
$ws.Range("D99").Value = 10063
$ws.Range("E99").Value = 3992
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 213

# Row 111: Luxemburgo
$ws.Range("B111").Value = 9840
$ws.Range("C111").Value = 109
$ws.Range("D111").Value = 8248
$ws.Range("E111").Value = 1459

# Row 116: Jamaica
$ws.Range("B116").Value = 7910
$ws.Range("C116").Value = 97
$ws.Range("D116").Value = 3303
$ws.Range("E116").Value = 4461

# Row 122: Cuba
$ws.Range("B122").Value = 6017
$ws.Range("C122").Value = 17
$ws.Range("D122").Value = 5602
$ws.Range("E122").Value = 292

# Row 134: Sri Lanka
$ws.Range("B134").Value = 4893
$ws.Range("C134").Value = 49
$ws.Range("D134").Value = 3328
$ws.Range("E134").Value = 1552
$ws.Range("H134").Value = 13

# Row 135: Republica de Africa Central
$ws.Range("B135").Value = 4854
$ws.Range("D135").Value = 1924
$ws.Range("E135").Value = 2868
$ws.Range("H135").Value = 62

# "Last refreshed" timestamp label
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 17:14"
